# feat(core): apply global row rules to Schema Editor and update test data
# Insert a new dummy row at the very top of the sheet so the 1st row is a
# dummy row (ignored by the binder) and the 2nd row becomes the header row,
# matching the "global row rules" used elsewhere. This pushes all existing
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything down by inserting a new row at the top.
$ws.Rows.Item(1).Insert()

# Populate the new first row with the dummy-row marker text (only column A).
$ws.Range("A1").Value = "# Dummy Row (Ignored by ExcelBinder)"
